$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (20) for hour 19 by copying the format of the last
# existing data row (19), then overwrite with the new data values.
$ws.Range("A19:T19").Copy()
$ws.Range("A20:T20").PasteSpecial()

$row = 20
$ws.Range("B$row").Value = 19
$ws.Range("C$row").Value = 0
$ws.Range("D$row").Value = 0
$ws.Range("E$row").Value = 0
$ws.Range("F$row").Value = 15680
$ws.Range("G$row").Value = 0
$ws.Range("H$row").Value = 0
$ws.Range("I$row").Value = 0
$ws.Range("J$row").Value = 0
$ws.Range("K$row").Value = 15680
$ws.Range("L$row").Value = 0
$ws.Range("M$row").Value = 0
$ws.Range("N$row").Value = 0
$ws.Range("O$row").Value = 15681.568
$ws.Range("P$row").Value = 0
$ws.Range("Q$row").Value = 0
$ws.Range("R$row").Value = 0
$ws.Range("S$row").Value = 0
$ws.Range("T$row").Value = 15681.568
